$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new transformer parameter headers
$ws.Range("K1").Value = "d_model"
$ws.Range("L1").Value = "num_layers"

# Add their corresponding values
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 2

# Update selection to match the saved workbook state
$ws.Range("J5").Select()
